$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Overview" (File Name / Path And Name / Extension / Publish URL /
# zh-cn / de-de / Latest HO Xliff Generate Date) gains 3 rows for the newly
# handed-off markdown file and its two image dependencies.
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$loOverview = $wsOverview.ListObjects.Item(1)

$wsOverview.Rows.Item(4).Copy()
$wsOverview.Rows.Item(5).Insert()
$wsOverview.Rows.Item(4).Copy()
$wsOverview.Rows.Item(6).Insert()
$wsOverview.Rows.Item(4).Copy()
$wsOverview.Rows.Item(7).Insert()
$loOverview.Resize($wsOverview.Range("A1:G7"))

$wsOverview.Range("A5").Value = "eb8ab756-3dbf-4572-8f49-25b5607fafd0.md"
$wsOverview.Range("C5").Value = ".md"
$wsOverview.Range("D5").Value = ""
$wsOverview.Range("E5").Value = "Ready for handoff"
$wsOverview.Range("F5").Value = "Ready for handoff"
$wsOverview.Range("G5").Value = "2016-11-02 05:50:11"

$wsOverview.Range("A6").Value = "f60228d2-3616-46ba-b922-79e8c510cab1.png"
$wsOverview.Range("C6").Value = ".png"
$wsOverview.Range("D6").Value = ""
$wsOverview.Range("E6").Value = "Ready for handoff"
$wsOverview.Range("F6").Value = "Ready for handoff"
$wsOverview.Range("G6").Value = "2016-11-02 05:50:11"

$wsOverview.Range("A7").Value = "ec5f1808-f645-4f30-b119-7c1fe8985a7e.png"
$wsOverview.Range("C7").Value = ".png"
$wsOverview.Range("D7").Value = ""
$wsOverview.Range("E7").Value = "Ready for handoff"
$wsOverview.Range("F7").Value = "Ready for handoff"
$wsOverview.Range("G7").Value = "2016-11-02 05:50:11"

$wsOverview.Hyperlinks.Add($wsOverview.Range("B5"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/50620f87aaf671ce5a8a5d4070077692ffc39335/e2e/eb8ab756-3dbf-4572-8f49-25b5607fafd0.md", "", "", "e2e\eb8ab756-3dbf-4572-8f49-25b5607fafd0.md")
$wsOverview.Hyperlinks.Add($wsOverview.Range("B6"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/50620f87aaf671ce5a8a5d4070077692ffc39335/e2e/f60228d2-3616-46ba-b922-79e8c510cab1.png", "", "", "e2e\f60228d2-3616-46ba-b922-79e8c510cab1.png")
$wsOverview.Hyperlinks.Add($wsOverview.Range("B7"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/50620f87aaf671ce5a8a5d4070077692ffc39335/e2e/ec5f1808-f645-4f30-b119-7c1fe8985a7e.png", "", "", "e2e\ec5f1808-f645-4f30-b119-7c1fe8985a7e.png")

# ---------------------------------------------------------------------------
# Sheet "zh-cn" gains the matching detail rows (table1.xml).
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")
$loZh = $wsZh.ListObjects.Item(1)

$wsZh.Rows.Item(4).Copy()
$wsZh.Rows.Item(5).Insert()
$wsZh.Rows.Item(4).Copy()
$wsZh.Rows.Item(6).Insert()
$wsZh.Rows.Item(4).Copy()
$wsZh.Rows.Item(7).Insert()
$loZh.Resize($wsZh.Range("A1:P7"))

$wsZh.Range("B5").Value = ".md"
$wsZh.Range("C5").Value = "Ready for handoff"
$wsZh.Range("D5").Value = "e2e"
$wsZh.Range("E5").Value = "ht"
$wsZh.Range("F5").Value = "False"
$wsZh.Range("G5").Value = "eb8ab756-3dbf-4572-8f49-25b5607fafd0.f29cb93078b07eed539eb380cbdb87114bf032ab.zh-cn.xlf"
$wsZh.Range("H5").Value = "2016-11-02 05:49:54"
$wsZh.Range("I5").Style = "Normal"
$wsZh.Range("I5").Value = ""
$wsZh.Range("J5").Value = ""
$wsZh.Range("K5").Value = "0001-01-01 00:00:00"
$wsZh.Range("L5").Value = ""
$wsZh.Range("M5").Value = "True"
$wsZh.Range("N5").Value = ""
$wsZh.Range("O5").Value = "False"
$wsZh.Range("P5").Value = ""

$wsZh.Range("B6").Value = ".png"
$wsZh.Range("C6").Value = "Ready for handoff"
$wsZh.Range("D6").Value = "e2e"
$wsZh.Range("E6").Value = "ht"
$wsZh.Range("F6").Value = "False"
$wsZh.Range("G6").Value = "2a546161cd49519b90afbd247352939c0ff0198e.png"
$wsZh.Range("H6").Value = "2016-11-02 05:49:54"
$wsZh.Range("I6").Style = "Normal"
$wsZh.Range("I6").Value = ""
$wsZh.Range("J6").Value = ""
$wsZh.Range("K6").Value = "0001-01-01 00:00:00"
$wsZh.Range("L6").Value = ""
$wsZh.Range("M6").Value = "True(Dependency)"
$wsZh.Range("N6").Value = "e2e\eb8ab756-3dbf-4572-8f49-25b5607fafd0.md"
$wsZh.Range("O6").Value = "False"
$wsZh.Range("P6").Value = ""

$wsZh.Range("B7").Value = ".png"
$wsZh.Range("C7").Value = "Ready for handoff"
$wsZh.Range("D7").Value = "e2e"
$wsZh.Range("E7").Value = "ht"
$wsZh.Range("F7").Value = "False"
$wsZh.Range("G7").Value = "d08688717f888f3f4950ca95b1021c8310d4b96e.png"
$wsZh.Range("H7").Value = "2016-11-02 05:49:54"
$wsZh.Range("I7").Style = "Normal"
$wsZh.Range("I7").Value = ""
$wsZh.Range("J7").Value = ""
$wsZh.Range("K7").Value = "0001-01-01 00:00:00"
$wsZh.Range("L7").Value = ""
$wsZh.Range("M7").Value = "True(Dependency)"
$wsZh.Range("N7").Value = "e2e\eb8ab756-3dbf-4572-8f49-25b5607fafd0.md"
$wsZh.Range("O7").Value = "False"
$wsZh.Range("P7").Value = ""

$wsZh.Hyperlinks.Add($wsZh.Range("A5"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/05fbcdff76596c48403f27016291490fa0299b52/e2e/eb8ab756-3dbf-4572-8f49-25b5607fafd0.md", "", "", "eb8ab756-3dbf-4572-8f49-25b5607fafd0.md")
$wsZh.Hyperlinks.Add($wsZh.Range("A6"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/05fbcdff76596c48403f27016291490fa0299b52/e2e/f60228d2-3616-46ba-b922-79e8c510cab1.png", "", "", "f60228d2-3616-46ba-b922-79e8c510cab1.png")
$wsZh.Hyperlinks.Add($wsZh.Range("A7"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/05fbcdff76596c48403f27016291490fa0299b52/e2e/ec5f1808-f645-4f30-b119-7c1fe8985a7e.png", "", "", "ec5f1808-f645-4f30-b119-7c1fe8985a7e.png")

$wsZh.Columns.Item(13).ColumnWidth = 16.83
$wsZh.Columns.Item(14).ColumnWidth = 39.15

# ---------------------------------------------------------------------------
# Sheet "de-de" mirrors "zh-cn" (table2.xml).
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")
$loDe = $wsDe.ListObjects.Item(1)

$wsDe.Rows.Item(4).Copy()
$wsDe.Rows.Item(5).Insert()
$wsDe.Rows.Item(4).Copy()
$wsDe.Rows.Item(6).Insert()
$wsDe.Rows.Item(4).Copy()
$wsDe.Rows.Item(7).Insert()
$loDe.Resize($wsDe.Range("A1:P7"))

$wsDe.Range("B5").Value = ".md"
$wsDe.Range("C5").Value = "Ready for handoff"
$wsDe.Range("D5").Value = "e2e"
$wsDe.Range("E5").Value = "ht"
$wsDe.Range("F5").Value = "False"
$wsDe.Range("G5").Value = "eb8ab756-3dbf-4572-8f49-25b5607fafd0.f29cb93078b07eed539eb380cbdb87114bf032ab.de-de.xlf"
$wsDe.Range("H5").Value = "2016-11-02 05:50:11"
$wsDe.Range("I5").Style = "Normal"
$wsDe.Range("I5").Value = ""
$wsDe.Range("J5").Value = ""
$wsDe.Range("K5").Value = "0001-01-01 00:00:00"
$wsDe.Range("L5").Value = ""
$wsDe.Range("M5").Value = "True"
$wsDe.Range("N5").Value = ""
$wsDe.Range("O5").Value = "False"
$wsDe.Range("P5").Value = ""

$wsDe.Range("B6").Value = ".png"
$wsDe.Range("C6").Value = "Ready for handoff"
$wsDe.Range("D6").Value = "e2e"
$wsDe.Range("E6").Value = "ht"
$wsDe.Range("F6").Value = "False"
$wsDe.Range("G6").Value = "2a546161cd49519b90afbd247352939c0ff0198e.png"
$wsDe.Range("H6").Value = "2016-11-02 05:50:11"
$wsDe.Range("I6").Style = "Normal"
$wsDe.Range("I6").Value = ""
$wsDe.Range("J6").Value = ""
$wsDe.Range("K6").Value = "0001-01-01 00:00:00"
$wsDe.Range("L6").Value = ""
$wsDe.Range("M6").Value = "True(Dependency)"
$wsDe.Range("N6").Value = "e2e\eb8ab756-3dbf-4572-8f49-25b5607fafd0.md"
$wsDe.Range("O6").Value = "False"
$wsDe.Range("P6").Value = ""

$wsDe.Range("B7").Value = ".png"
$wsDe.Range("C7").Value = "Ready for handoff"
$wsDe.Range("D7").Value = "e2e"
$wsDe.Range("E7").Value = "ht"
$wsDe.Range("F7").Value = "False"
$wsDe.Range("G7").Value = "d08688717f888f3f4950ca95b1021c8310d4b96e.png"
$wsDe.Range("H7").Value = "2016-11-02 05:50:11"
$wsDe.Range("I7").Style = "Normal"
$wsDe.Range("I7").Value = ""
$wsDe.Range("J7").Value = ""
$wsDe.Range("K7").Value = "0001-01-01 00:00:00"
$wsDe.Range("L7").Value = ""
$wsDe.Range("M7").Value = "True(Dependency)"
$wsDe.Range("N7").Value = "e2e\eb8ab756-3dbf-4572-8f49-25b5607fafd0.md"
$wsDe.Range("O7").Value = "False"
$wsDe.Range("P7").Value = ""

$wsDe.Hyperlinks.Add($wsDe.Range("A5"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/666f5f1156e6c96816c52ce400b4e901d4a24747/e2e/eb8ab756-3dbf-4572-8f49-25b5607fafd0.md", "", "", "eb8ab756-3dbf-4572-8f49-25b5607fafd0.md")
$wsDe.Hyperlinks.Add($wsDe.Range("A6"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/666f5f1156e6c96816c52ce400b4e901d4a24747/e2e/f60228d2-3616-46ba-b922-79e8c510cab1.png", "", "", "f60228d2-3616-46ba-b922-79e8c510cab1.png")
$wsDe.Hyperlinks.Add($wsDe.Range("A7"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/666f5f1156e6c96816c52ce400b4e901d4a24747/e2e/ec5f1808-f645-4f30-b119-7c1fe8985a7e.png", "", "", "ec5f1808-f645-4f30-b119-7c1fe8985a7e.png")

$wsDe.Columns.Item(13).ColumnWidth = 16.83
$wsDe.Columns.Item(14).ColumnWidth = 39.15
